# 8.10.1 indicator sheet: add the 2020 data column (Q) after the existing
# 2007-2019 columns (D..P), mirroring the year-over-year layout already in
# the sheet, then leave the selection on the newly-entered figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New 2020 column data -------------------------------------------------
# Year header
$ws.Range("Q3").Value = 2020

# Raw source figures (commercial bank branches, ATMs, adult population)
$ws.Range("Q6").Value = 312
$ws.Range("Q7").Value = 1856
$ws.Range("Q8").Value = 4337617

# Derived per-100,000-adults ratios, same formula pattern as column P
$ws.Range("Q4").Formula = "=Q6/Q8*100000"
$ws.Range("Q5").Formula = "=Q7/Q8*100000"

# --- Match formatting of the preceding (2019 / column P) cells -----------
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)

$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)

$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)

$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Re-assert the values/formulas: PasteSpecial(formats) must not clobber them
$ws.Range("Q3").Value = 2020
$ws.Range("Q6").Value = 312
$ws.Range("Q7").Value = 1856
$ws.Range("Q8").Value = 4337617
$ws.Range("Q4").Formula = "=Q6/Q8*100000"
$ws.Range("Q5").Formula = "=Q7/Q8*100000"

# --- Selection left on the newly added column, scrolled into view --------
$ws.Range("Q12").Select()
